# =====================================================================
#  evaluation_scores.xlsx -- "add lr = 0.05"
#  1) Add a new "Count of data" sheet (train/test class-count tables)
#  2) Add a "learning rate = 0.05" results block to the BOG sheet
#  3) Make BOG the active/selected sheet (was Comparison)
# =====================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) New sheet "Count of data", appended after "Comparison"
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$data = $wb.Worksheets.Add($null, $lastSheet)
$data.Name = "Count of data"

# -- section titles (row 3) --
$data.Range("A3").Value = "Traning Data"
$data.Range("F3").Value = "Test Data"

# -- column headers (row 5), highlighted yellow --
$data.Range("A5").Value = "Emotion"
$data.Range("B5").Value = "N"
$data.Range("C5").Value = "Percentage"
$data.Range("F5").Value = "Emotion"
$data.Range("G5").Value = "N"
$data.Range("H5").Value = "Percentage"
$data.Range("A5:C5").Interior.Color = 65535
$data.Range("F5:H5").Interior.Color = 65535

# -- per-emotion counts, training (A:C) and test (F:H) --
$labels = "JOY","FEAR","SHAME","DISGUST","GUILT","ANGER","SADNESS"
$trainN = 777,751,758,758,768,758,761
$testN  = 162,164,164,173,155,176,151
$rows   = 6,7,8,9,10,11,12

for ($i = 0; $i -lt $labels.Length; $i++) {
    $r = $rows[$i]
    $data.Range("A$r").Value = $labels[$i]
    $data.Range("B$r").Value = $trainN[$i]
    $data.Range("C$r").Formula = "=B$r/B13"
    $data.Range("F$r").Value = $labels[$i]
    $data.Range("G$r").Value = $testN[$i]
    $data.Range("H$r").Formula = "=G$r/G13"
}
$data.Range("A6:B12").Font.Family = 1
$data.Range("F6:G12").Font.Family = 1

# -- totals row --
$data.Range("B13").Formula = "=SUM(B6:B12)"
$data.Range("G13").Formula = "=SUM(G6:G12)"
$data.Range("J13").Value = 6476

# -- footnotes --
$data.Range("B14").Value = "~=82.3%"
$data.Range("G14").Value = "~= 17.7%"

$data.Range("H10").Select()

# ---------------------------------------------------------------------
# 2) BOG sheet: new "learning rate = 0.05" block (rows 18-32)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("BOG")

# -- banner row --
$ws.Range("A18").Value = "learning rate = 0.05"
$ws.Range("A18").Font.ThemeColor = 1
$ws.Range("A18").Interior.Color = 255
$ws.Range("B18").Font.ThemeColor = 1

# -- confusion-matrix text, entered in the same (slightly out-of-order)
#    sequence as the original author --
$ws.Range("B31").Value = "[[ 75 127]"
$ws.Range("B32").Value = " [ 80 864]] "

$ws.Range("B19").Value = "[[102  59]"
$ws.Range("B20").Value = " [ 60 925]] "

$ws.Range("B21").Value = "[[ 49  78]"
$ws.Range("B22").Value = " [127 892]]"

$ws.Range("B23").Value = "[[ 54  51]"
$ws.Range("B24").Value = " [110 931]] "

$ws.Range("B25").Value = "[[ 95  95]"
$ws.Range("B26").Value = " [ 78 878]]"

$ws.Range("B27").Value = "[[ 85  80]"
$ws.Range("B28").Value = " [ 66 915]] "

$ws.Range("B30").Value = " [ 69 940]]"
$ws.Range("B29").Value = "[[ 95  42]"

# -- emotion labels + F-scores --
$ws.Range("A19").Value = "JOY"
$ws.Range("C19").Value = 0.63157894736842102

$ws.Range("A21").Value = "ANGER"
$ws.Range("C21").Value = 0.32343234323432302

$ws.Range("A23").Value = "SHAME"
$ws.Range("C23").Value = 0.40148698884758299

$ws.Range("A25").Value = "DISGUST"
$ws.Range("C25").Value = 0.52341597796143202

$ws.Range("A27").Value = "SADNESS"
$ws.Range("C27").Value = 0.537974683544303

$ws.Range("A29").Value = "FEAR"
$ws.Range("C29").Value = 0.63122923588039803

$ws.Range("A31").Value = "GUILT"
$ws.Range("C31").Value = 0.42016806722688999

# -- wrap text on the two rows that need the extra height --
$ws.Range("B21").WrapText = $true
$ws.Range("B21").RowHeight = 16
$ws.Range("B29").WrapText = $true
$ws.Range("B29").RowHeight = 16

# -- widen column A now that it holds longer labels --
$ws.Range("A1").ColumnWidth = 17.285714285714285

# ---------------------------------------------------------------------
# 3) BOG becomes the active / selected sheet
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("B18").Select()
